$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "60.940.68"
Set-TextValue "E2" "  -2.26%  "
Set-TextValue "D3" "3.377.43"
Set-TextValue "E3" "  -1.13%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "572.47"
Set-TextValue "E5" "  -1.89%  "
Set-TextValue "D6" "136.04"
Set-TextValue "E6" "  -2.74%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "3.376.79"
Set-TextValue "E8" "  -1.08%  "
Set-TextValue "D9" "0.469"
Set-TextValue "E9" "  -1.88%  "
Set-TextValue "E10" "  +1.24%  "
Set-TextValue "D11" "0.122"
Set-TextValue "E11" "  -4.77%  "
Set-TextValue "D12" "0.383"
Set-TextValue "E12" "  -2.92%  "
Set-TextValue "D13" "3.955.23"
Set-TextValue "E13" "  -1.09%  "
Set-TextValue "E14" "  +0.76%  "
Set-TextValue "B15" "ShibaInu"
Set-TextValue "C15" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.0000172"
Set-TextValue "E15" "  -4.81%  "
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.373.45"
Set-TextValue "E16" "  -1.42%  "
Set-TextValue "E17" "  -0.95%  "
Set-TextValue "D18" "61.113.32"
Set-TextValue "E18" "  -1.99%  "
Set-TextValue "D19" "13.82"
Set-TextValue "E19" "  -3.10%  "
Set-TextValue "D20" "5.75"
Set-TextValue "E20" "  -1.93%  "
Set-TextValue "D21" "9.31"
Set-TextValue "E21" "  -3.40%  "
Set-TextValue "D22" "374.40"
Set-TextValue "E22" "  -5.44%  "
Set-TextValue "D23" "3.513.99"
Set-TextValue "E23" "  -0.96%  "
Set-TextValue "D24" "0.549"
Set-TextValue "E24" "  -3.84%  "
Set-TextValue "E25" "  +0.03%  "
Set-TextValue "D26" "0.0000125"
Set-TextValue "E26" "  -5.47%  "
Set-TextValue "D27" "70.93"
Set-TextValue "E27" "  -1.33%  "
Set-TextValue "B28" "Kaspa"
Set-TextValue "C28" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.181"
Set-TextValue "E28" "  +10.85%  "
Set-TextValue "B29" "Fetch.AI"
Set-TextValue "C29" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D29" "1.67"
Set-TextValue "E29" "  +2.28%  "
Set-TextValue "B30" "Binance-PegBSC-USD"
Set-TextValue "C30" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.05%  "
Set-TextValue "B31" "RenderToken"
Set-TextValue "C31" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D31" "7.39"
Set-TextValue "E31" "  -5.15%  "
Set-TextValue "D32" "8.06"
Set-TextValue "E32" "  -3.61%  "
Set-TextValue "E33" "  -2.86%  "
Set-TextValue "D35" "23.49"
Set-TextValue "E35" "  -0.73%  "
Set-TextValue "D36" "5.19"
Set-TextValue "E36" "  -5.69%  "
Set-TextValue "D37" "1.54"
Set-TextValue "E37" "  -3.59%  "
Set-TextValue "D38" "6.80"
Set-TextValue "E38" "  -2.82%  "
Set-TextValue "D39" "164.55"
Set-TextValue "E39" "  -0.11%  "
Set-TextValue "D40" "0.0752"
Set-TextValue "E40" "  -5.69%  "
Set-TextValue "E41" "  +0.04%  "
Set-TextValue "D42" "0.773"
Set-TextValue "E42" "  -1.79%  "
Set-TextValue "B43" "OKB"
Set-TextValue "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "41.75"
Set-TextValue "E43" "  -0.47%  "
Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "25.00"
Set-TextValue "E44" "  -1.06%  "
Set-TextValue "B45" "Stacks"
Set-TextValue "C45" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "1.69"
Set-TextValue "E45" "  -6.76%  "
Set-TextValue "D46" "1.19"
Set-TextValue "E46" "  -7.00%  "
Set-TextValue "D47" "4.35"
Set-TextValue "E47" "  -3.20%  "
Set-TextValue "D48" "2.462.90"
Set-TextValue "E48" "  +4.06%  "
Set-TextValue "B49" "Cosmos"
Set-TextValue "C49" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "6.78"
Set-TextValue "E49" "  -2.70%  "
Set-TextValue "B50" "InjectiveProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D50" "22.89"
Set-TextValue "E50" "  -2.75%  "
Set-TextValue "D51" "2.41"
Set-TextValue "E51" "  +2.33%  "
